$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D reference-range text fixes (replace ~ with -, strip units) ---
$ws.Range("D2").Value = "9.0-50.0"
$ws.Range("D3").Value = "15.0-40.0"
$ws.Range("D4").Value = "5.0-22.0"
$ws.Range("D5").Value = "1.7-7.1"
$ws.Range("D6").Value = "65.0-85.0"
$ws.Range("D7").Value = "40.0-55.0"
$ws.Range("D8").Value = "20.00-30.00"
$ws.Range("D9").Value = "1.20-2.40"
$ws.Range("D10").Value = "45.0-125.0"
$ws.Range("D11").Value = "3.9-6.1"
$ws.Range("D12").Value = "3.00-6.00"
$ws.Range("D13").Value = "340.50-1.70"
$ws.Range("D14").Value = "0.80-2.00"
$ws.Range("D15").Value = "2.1-3.1"

# --- Column B label text fixes (typo corrections) ---
$ws.Range("B4").Value = "总胆红素"
$ws.Range("B9").Value = "载脂蛋白比值"
$ws.Range("B13").Value = "甘油三酯"
$ws.Range("B14").Value = "高密度脂蛋白胆固醇"
$ws.Range("B15").Value = "低密度脂蛋白胆固醇"

# --- Column C result value fixes (force text format so numeric-looking
#     strings like "31.2" stay text, matching the source inlineStr cells) ---
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "31.2"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "1.44"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "1.8"

# --- Column A: clear duplicated reference-range junk rows 16-22, 24-29 ---
# (touch NumberFormat first so the row/cell is preserved as a blank entry
#  instead of being dropped entirely when emptied)
$ws.Range("A16:A22").NumberFormat = "@"
$ws.Range("A16:A22").Value = ""
$ws.Range("A24:A29").NumberFormat = "@"
$ws.Range("A24:A29").Value = ""

$ws.Range("A23").Value = ".20~2."
